$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "EmpID" header in column E
$ws.Range("E1").Value = "EmpID"

# Make the header row bold (A1:E1), matching the Normal 11pt base font
$ws.Range("A1:E1").Font.Size = 11
$ws.Range("A1:E1").Font.Bold = $true

# Set header row height
$ws.Rows.Item(1).RowHeight = 15

# Update selection to match the target state
$ws.Range("A1:E1").Select()
